$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.073.10'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '4.009.89'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.696'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +10.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.746'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('E11').Value = '  -5.81%  '
$ws.Range('E12').Value = '  +5.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.67'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').Value = '4.648.30'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '4.002.79'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.00'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.63%  '
$ws.Range('D20').Value = '71.902.17'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '426.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '97.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.27%  '
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.77'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.61'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +24.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.129'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.11'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '673.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.53%  '
$ws.Range('E35').Value = '  -3.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '43.15'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.427'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.153'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('E39').Value = '  -8.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.42'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.13%  '
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('E45').Value = '  +2.27%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.80%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.78%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.63'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.02%  '
$ws.Range('E49').Value = '  -7.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000271'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.71'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.57%  '
